$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new resistor (R27) joined the zero-ohm group on the BOM, so the
# quantity and reference designators for that row need updating.
$ws.Range("A3").Value = 3
$ws.Range("E3").Value = "R6, R27, R28"

# Final release of V1.0.0 adds a new part: SW3, a DPDT slide toggle switch.
$ws.Range("A32").Value = 1
$ws.Range("B32").Value = "JS202011JCQN"
$ws.Range("C32").Value = "JS202011JCQN"
$ws.Range("D32").Value = "JS202011JCQN"
$ws.Range("E32").Value = "SW3"
$ws.Range("F32").Value = "SLIDE TOGGLE SWITCH VERTICAL (DPDT)"
$ws.Range("G32").Value = "C&K Components"
$ws.Range("H32").Value = "JS202011JCQN"
$ws.Range("I32").Value = "CKN10723CT-ND"
$ws.Range("J32").Value = "611-JS202011JCQN"

# Leave the cursor where the author last left it when saving.
$ws.Range("E41").Select() | Out-Null
